# The underlying data source was re-synced and the rows for this sighting
# batch (Id 4-12) come back in a different order. Re-map each row to the
# record it now corresponds to by snapshotting the "before" values first
# (so overwrites don't clobber a value we still need to read), then writing
# them into their new homes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Get-RowData($r) {
    $data = @{
        A = $ws.Cells.Item($r, 1).Value2
        B = $ws.Cells.Item($r, 2).Value2
        D = $ws.Cells.Item($r, 4).Value2
        E = $ws.Cells.Item($r, 5).Value2
        F = $ws.Cells.Item($r, 6).Value2
        G = $ws.Cells.Item($r, 7).Value2
        H = $ws.Cells.Item($r, 8).Value2
        Q = $ws.Cells.Item($r, 17).Value2
        R = $ws.Cells.Item($r, 18).Value2
    }
    return $data
}

function Set-RowData($r, $data) {
    $ws.Cells.Item($r, 1).Value = $data.A
    $ws.Cells.Item($r, 2).Value = $data.B
    $ws.Cells.Item($r, 4).Value = $data.D
    $ws.Cells.Item($r, 5).Value = $data.E
    $ws.Cells.Item($r, 6).Value = $data.F
    $ws.Cells.Item($r, 7).Value = $data.G
    $ws.Cells.Item($r, 8).Value = $data.H
    $ws.Cells.Item($r, 17).Value = $data.Q
    $ws.Cells.Item($r, 18).Value = $data.R
}

# Snapshot every source row (4-12) before any writes happen.
$row4 = Get-RowData 4
$row5 = Get-RowData 5
$row6 = Get-RowData 6
$row7 = Get-RowData 7
$row8 = Get-RowData 8
$row9 = Get-RowData 9
$row10 = Get-RowData 10
$row11 = Get-RowData 11
$row12 = Get-RowData 12

# Row 10 also carries bird-specific columns (K,L,M,N empty markers + AC
# "hack" comment) that must move to row 11 along with the rest of its data.
$row10K = $ws.Cells.Item(10, 11).Value2
$row10L = $ws.Cells.Item(10, 12).Value2
$row10M = $ws.Cells.Item(10, 13).Value2
$row10N = $ws.Cells.Item(10, 14).Value2
$row10AC = $ws.Cells.Item(10, 29).Value2

# New row -> source (old) row mapping.
Set-RowData 4 $row11
Set-RowData 5 $row12
Set-RowData 6 $row8
Set-RowData 7 $row5
Set-RowData 8 $row4
Set-RowData 9 $row9
Set-RowData 10 $row6
Set-RowData 11 $row10
Set-RowData 12 $row7

# Row 10 (now holding the old row 6/"Garnlav" record) no longer has the
# bird-only columns - clear them out.
$ws.Cells.Item(10, 11).Value = ""
$ws.Cells.Item(10, 12).Value = ""
$ws.Cells.Item(10, 13).Value = ""
$ws.Cells.Item(10, 14).Value = ""
$ws.Cells.Item(10, 29).Value = ""

# Row 11 (now holding the old row 10/"Spillkråka" record) picks up those
# bird-only columns instead.
$ws.Cells.Item(11, 11).Value = $row10K
$ws.Cells.Item(11, 12).Value = $row10L
$ws.Cells.Item(11, 13).Value = $row10M
$ws.Cells.Item(11, 14).Value = $row10N
$ws.Cells.Item(11, 29).Value = $row10AC
